$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 3) mirroring row 2's data (both "will")
$ws.Range("A3").Value = "will"
$ws.Range("B3").Value = "will"

# Select B3 like the final saved state shows
$ws.Range("B3").Select()
